$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -------------------------------------------------
# "preRequisites" goes in before the old column D (HttpMethod), pushing
# HttpMethod/Uri/Body/ExpectedStatusCode one column to the right.
$ws.Columns.Item(4).Insert()

# "FieldValidations" goes in right before ExpectedResponseBody, which (after the
# previous insert) now sits in column I.
$ws.Columns.Item(9).Insert()

# --- Row 2 (first, pre-existing test case) --------------------------------------
$ws.Range("F2").Value = "impacts/v1/impacts?startDate=2020-5-1&endDate=2022-5-1"
$ws.Range("H2").Value = "400"

# --- Row 3 (second, new successful test case) -----------------------------------
$ws.Range("F3").Value = "impacts/v1/impacts?startDate=2020-5-1&endDate=2022-5-1"

# --- Header row (row 1) ---------------------------------------------------------
$ws.Range("A1").Value = "TCID"
$ws.Range("B1").Value = "Run"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "preRequisites"
$ws.Range("E1").Value = "HttpMethod"
$ws.Range("F1").Value = "Uri"
$ws.Range("G1").Value = "Body"
$ws.Range("H1").Value = "ExpectedStatusCode"
$ws.Range("I1").Value = "FieldValidations"
$ws.Range("J1").Value = "ExpectedResponseBody"

# --- Drop the blank placeholder cells the column-insert left behind -------------
$ws.Range("D2:D3").Clear()

# --- Column widths ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 13.833333333333334
$ws.Columns.Item(5).ColumnWidth = 11
$ws.Columns.Item(6).ColumnWidth = 49.5
$ws.Columns.Item(7).ColumnWidth = 20.833333333333332
$ws.Columns.Item(8).ColumnWidth = 16.666666666666668
$ws.Columns.Item(9).ColumnWidth = 16.666666666666668
$ws.Columns.Item(10).ColumnWidth = 25.333333333333332

# --- View / selection -------------------------------------------------------------
$ws.Range("G5").Select()
